# New Submission Synced: 2026-02-08 22:15:58
# Sheet "JSS 3F" is the Google-Forms-style response log. A previous row
# (row 6) had its "Admission No" typed as the text "05" instead of the
# numeric 5 used everywhere else in that column; this sync also normalizes
# it to the number 5. A brand-new response (submitted 2026-02-08 22:15:58)
# is appended as row 7, again with the admission number typed as "05" and
# therefore stored as text to preserve the leading zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3F")

# Normalize C6 ("Admission No") from text "05" to the number 5.
$ws.Range("C6").Value = 5

# Append the new submission as row 7.
$ws.Range("A7").Value = "2026-02-08 22:15:58"
$ws.Range("B7").Value = "Usman Muhammad Gubio"
# Leading apostrophe forces text storage so the leading zero is kept,
# matching how "05" is stored elsewhere in this column.
$ws.Range("C7").Value = "'05"
$ws.Range("D7").Value = 9
